$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B width to match column A's stored width (15.42578125)
$ws.Columns("B").ColumnWidth = 14.6

# Update values in A1:B32
$ws.Range("A1").Value = -0.24806121495476674
$ws.Range("B1").Value = 0.24785262924747542
$ws.Range("A2").Value = -0.14829472432063007
$ws.Range("B2").Value = 0.14794202502065268
$ws.Range("A3").Value = -0.098240167336903994
$ws.Range("B3").Value = 0.097941941583554737
$ws.Range("A4").Value = -0.089941941648937984
$ws.Range("B4").Value = 0.089545088584159771
$ws.Range("A5").Value = -0.086545088621147848
$ws.Range("B5").Value = 0.085201537157288065
$ws.Range("A6").Value = -0.017723483871954215
$ws.Range("B6").Value = 0.017611866193286474
$ws.Range("A7").Value = -0.0076118662848196905
$ws.Range("B7").Value = 0.0075978055841514625
$ws.Range("A8").Value = 0.0024021943237504395
$ws.Range("B8").Value = -0.0024068553961833494
$ws.Range("A9").Value = 0.0044068553554752476
$ws.Range("B9").Value = -0.0044082616744054448
$ws.Range("A10").Value = 0.0064082616341121224
$ws.Range("B10").Value = -0.0064079191704937699
$ws.Range("A11").Value = 0.0094079191238103377
$ws.Range("B11").Value = -0.0094081515334130117
$ws.Range("A12").Value = 0.012908151483875141
$ws.Range("B12").Value = -0.012932782890636929
$ws.Range("A13").Value = -0.014937040943279634
$ws.Range("B13").Value = 0.014924523508414111
$ws.Range("A14").Value = -0.0069245235855062148
$ws.Range("B14").Value = 0.0069229352323798565
$ws.Range("A15").Value = -0.0059229352651311018
$ws.Range("B15").Value = 0.0059182916631099758
$ws.Range("A16").Value = -0.0060326376767783074
$ws.Range("B16").Value = 0.0060032012797917034
$ws.Range("A17").Value = -0.0040032013202990768
$ws.Range("B17").Value = 0.0039999999464104263
$ws.Range("A18").Value = -0.093557904967578764
$ws.Range("B18").Value = 0.093336294791125596
$ws.Range("A19").Value = -0.044024460439790136
$ws.Range("B19").Value = 0.043432623573657736
$ws.Range("A20").Value = -0.03943262361136135
$ws.Range("B20").Value = 0.039271178207499702
$ws.Range("A21").Value = -0.035271178247092472
$ws.Range("B21").Value = 0.035032666332921458
$ws.Range("A22").Value = -0.045701857719766537
$ws.Range("B22").Value = 0.045491676327090502
$ws.Range("A23").Value = -0.040491676371488872
$ws.Range("B23").Value = 0.040097515415445706
$ws.Range("A24").Value = -0.020097515559854173
$ws.Range("B24").Value = 0.019999999853611783
$ws.Range("A25").Value = -0.065478053325215058
$ws.Range("B25").Value = 0.065393994281212287
$ws.Range("A26").Value = -0.094569092521224363
$ws.Range("B26").Value = 0.09441305756401519
$ws.Range("A27").Value = -0.091913057611376026
$ws.Range("B27").Value = 0.090986654070008832
$ws.Range("A28").Value = -0.088986654127235276
$ws.Range("B28").Value = 0.088352468652152183
$ws.Range("A29").Value = -0.081352468751084928
$ws.Range("B29").Value = 0.081168767690726007
$ws.Range("A30").Value = -0.021168768133814186
$ws.Range("B30").Value = 0.021022604514183918
$ws.Range("A31").Value = -0.014022604620908652
$ws.Range("B31").Value = 0.014001007085269634
$ws.Range("A32").Value = -0.0040010072117091511
$ws.Range("B32").Value = 0.0039999999118336405
